# Regenerate the "K" column (column G) values for this save_data sheet.
# The commit updates the strike-count column (header "K", previously
# computed differently -> "K instead of Strike#") with freshly
# calculated values (std/mean regen + s_vals calc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column G ("K")
$newK = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 2
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 2
    15 = 1
    16 = 1
    17 = 1
    18 = 2
    19 = 2
    20 = 0
    21 = 1
    22 = 2
    23 = 0
    24 = 1
    27 = 1
    28 = 2
    29 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
